# Weekly update: insert two new price rows (week of 2022-01-17, serial 44578)
# for "Choclo" / "Choclero" and "Choclo" / "Dulce o Americano", quality
# "Segunda", at the top of the existing data block (rows 165-166), pushing
# all the previously-existing rows 165-175 down to 167-177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block.
$ws.Rows.Item(165).Insert()
$ws.Rows.Item(165).Insert()

# New row 165: Choclero / Segunda
$ws.Cells.Item(165, 1).Value = 4
$ws.Cells.Item(165, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(165, 3).Value = "Los Lagos"
$ws.Cells.Item(165, 4).Value = 44578
$ws.Cells.Item(165, 5).Value = 10
$ws.Cells.Item(165, 6).Value = 100112024
$ws.Cells.Item(165, 7).Value = "Choclo"
$ws.Cells.Item(165, 8).Value = "Choclero"
$ws.Cells.Item(165, 9).Value = "Segunda"
$ws.Cells.Item(165, 10).Value = 10000
$ws.Cells.Item(165, 11).Value = 350
$ws.Cells.Item(165, 12).Value = 350
$ws.Cells.Item(165, 13).Value = 350
$ws.Cells.Item(165, 14).Value = "`$/unidad"
$ws.Cells.Item(165, 15).Value = "Región del Maule"
$ws.Cells.Item(165, 16).Value = 350
$ws.Cells.Item(165, 17).Value = 1
$ws.Cells.Item(165, 18).Value = "Hortaliza"

# New row 166: Dulce o Americano / Segunda
$ws.Cells.Item(166, 1).Value = 4
$ws.Cells.Item(166, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(166, 3).Value = "Los Lagos"
$ws.Cells.Item(166, 4).Value = 44578
$ws.Cells.Item(166, 5).Value = 10
$ws.Cells.Item(166, 6).Value = 100112024
$ws.Cells.Item(166, 7).Value = "Choclo"
$ws.Cells.Item(166, 8).Value = "Dulce o Americano"
$ws.Cells.Item(166, 9).Value = "Segunda"
$ws.Cells.Item(166, 10).Value = 20000
$ws.Cells.Item(166, 11).Value = 200
$ws.Cells.Item(166, 12).Value = 200
$ws.Cells.Item(166, 13).Value = 200
$ws.Cells.Item(166, 14).Value = "`$/unidad"
$ws.Cells.Item(166, 15).Value = "Región del Maule"
$ws.Cells.Item(166, 16).Value = 200
$ws.Cells.Item(166, 17).Value = 1
$ws.Cells.Item(166, 18).Value = "Hortaliza"
